# Rename the two header blocks (old "_old"/"_new" suffixes -> per-input-file
# "_FV2210"/"_FV2304" suffixes) and turn the header row into a frozen,
# filterable Excel Table ("Table1") spanning the full used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J ("*_old") -> "*_FV2210"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2210"
}

# Column K ("diff") is unchanged.

# Columns L-U ("*_new") -> "*_FV2304"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2304"
}

# Freeze the header row (row 1).
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn A1:U66 into a native Excel Table with the new header names.
$tableRange = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
